$d = $word.ActiveDocument

# Paragraph 2: job title / location line
$p2 = $d.Paragraphs(2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$r2.Text = "Responsable de Projets Techniques depuis plus de 25 ans, je cherche une opportunité en tant que Product Owner ou Scrum Master près de Grenoble ou d’Annecy."

# Paragraph 3: autonomy / skills sentence
$p3 = $d.Paragraphs(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End)
$r3.Text = "Mon autonomie et mes connaissances des systèmes d'information, ainsi que mon expérience dans la gestion de projets et l’animation d’équipes me conduisent à penser que je serai un élément déterminant."

# Paragraph 4: agile experience sentence
$p4 = $d.Paragraphs(4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End)
$r4.Text = "Les sujets sur lesquels j'ai été amené à travailler, notamment en mode agile (Scrum), m'ont permis de développer une bonne maîtrise dans l’accompagnement d’équipes, l’organisation des projets et leur coordination."
